$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1: update submitted values with a fresh timestamp
$ws.Range("A1").Value = "adsfsf"
$ws.Range("B1").Value = "sdfsf"
$ws.Range("C1").Value = "Climate/Weather"
$ws.Range("D1").Value = 45189.46088765046

# Row 2: update submitted values with a fresh timestamp
$ws.Range("A2").Value = "eadfs"
$ws.Range("B2").Value = "sdfsfs"
$ws.Range("D2").Value = 45189.46135821888

# Clear selection: remove the stale row (error message row) entirely
$ws.Rows("3:3").Delete()

# Reflect the clear-selection click moving focus away from the grid
[void]$ws.Range("N17").Select()
